$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 12
$ws.Range("E2").Value = 11.27000045776367
$ws.Range("F2").Value = 13.59000015258789
$ws.Range("G2").Value = 10.89999961853027
$ws.Range("H2").Value = 651255398
$ws.Range("I2").Value = "CCCS"

$ws.Range("D3").Value = 12
$ws.Range("E3").Value = 11.27000045776367
$ws.Range("F3").Value = 13.59000015258789
$ws.Range("G3").Value = 10.89999961853027
$ws.Range("H3").Value = 651255398
$ws.Range("I3").Value = "CCCS"

$ws.Range("D4").Value = 12
$ws.Range("E4").Value = 11.27000045776367
$ws.Range("F4").Value = 13.59000015258789
$ws.Range("G4").Value = 10.89999961853027
$ws.Range("H4").Value = 651255398
$ws.Range("I4").Value = "CCCS"

$ws.Range("D5").Value = 12
$ws.Range("E5").Value = 11.27000045776367
$ws.Range("F5").Value = 13.59000015258789
$ws.Range("G5").Value = 10.89999961853027
$ws.Range("H5").Value = 651255398
$ws.Range("I5").Value = "CCCS"

$ws.Range("D6").Value = 13.25
$ws.Range("E6").Value = 13.47000026702881
$ws.Range("F6").Value = 15.47999954223633
$ws.Range("G6").Value = 12.75
$ws.Range("H6").Value = 651255398
$ws.Range("I6").Value = "CCCS"

$ws.Range("D7").Value = 10.14999961853027
$ws.Range("E7").Value = 10.26000022888184
$ws.Range("F7").Value = 10.40999984741211
$ws.Range("G7").Value = 10.0600004196167
$ws.Range("H7").Value = 651255398
$ws.Range("I7").Value = "CCCS"

$ws.Range("D8").Value = 9.960000038146973
$ws.Range("E8").Value = 9.260000228881836
$ws.Range("F8").Value = 10.01000022888184
$ws.Range("G8").Value = 8.760000228881836
$ws.Range("H8").Value = 651255398
$ws.Range("I8").Value = "CCCS"

$ws.Range("D9").Value = 10.51000022888184
$ws.Range("E9").Value = 11.97999954223633
$ws.Range("F9").Value = 12.27999973297119
$ws.Range("G9").Value = 9.899999618530272
$ws.Range("H9").Value = 651255398
$ws.Range("I9").Value = "CCCS"

$ws.Range("D10").Value = 11.35999965667725
$ws.Range("E10").Value = 10.65999984741211
$ws.Range("F10").Value = 11.42500019073486
$ws.Range("G10").Value = 9.630000114440918
$ws.Range("H10").Value = 651255398
$ws.Range("I10").Value = "CCCS"

$ws.Range("D11").Value = 11.01000022888184
$ws.Range("E11").Value = 9.229999542236328
$ws.Range("F11").Value = 11.53999996185303
$ws.Range("G11").Value = 9.170000076293944
$ws.Range("H11").Value = 651255398
$ws.Range("I11").Value = "CCCS"

$ws.Range("D12").Value = 9.199999809265137
$ws.Range("E12").Value = 9.989999771118164
$ws.Range("F12").Value = 10.03999996185303
$ws.Range("G12").Value = 8.909999847412109
$ws.Range("H12").Value = 651255398
$ws.Range("I12").Value = "CCCS"

$ws.Range("D13").Value = 9.140000343322754
$ws.Range("E13").Value = 9.329999923706056
$ws.Range("F13").Value = 9.520000457763672
$ws.Range("G13").Value = 8.194999694824219
$ws.Range("H13").Value = 651255398
$ws.Range("I13").Value = "CCCS"

$ws.Range("D14").Value = 8.840000152587891
$ws.Range("E14").Value = 9.25
$ws.Range("F14").Value = 9.340000152587891
$ws.Range("G14").Value = 8.300000190734863
$ws.Range("H14").Value = 651255398
$ws.Range("I14").Value = "CCCS"

$ws.Range("D15").Value = 8.970000267028809
$ws.Range("E15").Value = 8.680000305175781
$ws.Range("F15").Value = 9.079999923706056
$ws.Range("G15").Value = 8.539999961853027
$ws.Range("H15").Value = 651255398
$ws.Range("I15").Value = "CCCS"

$ws.Range("D16").Value = 11.10999965667725
$ws.Range("E16").Value = 11.02000045776367
$ws.Range("F16").Value = 11.21000003814697
$ws.Range("G16").Value = 10.36999988555908
$ws.Range("H16").Value = 651255398
$ws.Range("I16").Value = "CCCS"

$ws.Range("D17").Value = 13.35000038146973
$ws.Range("E17").Value = 10.77000045776367
$ws.Range("F17").Value = 13.35000038146973
$ws.Range("G17").Value = 10.56999969482422
$ws.Range("H17").Value = 651255398
$ws.Range("I17").Value = "CCCS"

$ws.Range("D18").Value = 11.34000015258789
$ws.Range("E18").Value = 10.98999977111816
$ws.Range("F18").Value = 11.35000038146973
$ws.Range("G18").Value = 10.60000038146973
$ws.Range("H18").Value = 651255398
$ws.Range("I18").Value = "CCCS"

$ws.Range("D19").Value = 11.97000026702881
$ws.Range("E19").Value = 11.22000026702881
$ws.Range("F19").Value = 12.05000019073486
$ws.Range("G19").Value = 11.19499969482422
$ws.Range("H19").Value = 651255398
$ws.Range("I19").Value = "CCCS"

$ws.Range("D20").Value = 11.13000011444092
$ws.Range("E20").Value = 10.26000022888184
$ws.Range("F20").Value = 11.51000022888184
$ws.Range("G20").Value = 10.05000019073486
$ws.Range("H20").Value = 651255398
$ws.Range("I20").Value = "CCCS"

$ws.Range("D21").Value = 11.05000019073486
$ws.Range("E21").Value = 10.40999984741211
$ws.Range("F21").Value = 11.39999961853027
$ws.Range("G21").Value = 10.38000011444092
$ws.Range("H21").Value = 651255398
$ws.Range("I21").Value = "CCCS"

$ws.Range("D22").Value = 11.80000019073486
$ws.Range("E22").Value = 11.10999965667725
$ws.Range("F22").Value = 11.80000019073486
$ws.Range("G22").Value = 10.9350004196167
$ws.Range("H22").Value = 651255398
$ws.Range("I22").Value = "CCCS"

$ws.Range("D23").Value = 9.029999732971191
$ws.Range("E23").Value = 9.260000228881836
$ws.Range("F23").Value = 9.350000381469728
$ws.Range("G23").Value = 8.140000343322754
$ws.Range("H23").Value = 651255398
$ws.Range("I23").Value = "CCCS"

$ws.Range("D24").Value = 9.380000114440918
$ws.Range("E24").Value = 9.670000076293944
$ws.Range("F24").Value = 10.07999992370606
$ws.Range("G24").Value = 9.164999961853027
$ws.Range("H24").Value = 651255398
$ws.Range("I24").Value = "CCCS"

